# Insert a new weekly price record as row 441 on the sheet, pushing the
# existing rows 441:460 down to 442:461 (dimension grows from R460 to R461).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(441).Insert()

$newRow = 441
$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 45041
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100112003
$ws.Cells.Item($newRow, 7).Value = "Ajo"
$ws.Cells.Item($newRow, 8).Value = "Chino"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 400
$ws.Cells.Item($newRow, 11).Value = 17000
$ws.Cells.Item($newRow, 12).Value = 18000
$ws.Cells.Item($newRow, 13).Value = 17500
$ws.Cells.Item($newRow, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item($newRow, 15).Value = "China"
$ws.Cells.Item($newRow, 16).Value = 1750
$ws.Cells.Item($newRow, 17).Value = 10
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date/time number format
# (style index "2" on column D for all data rows).
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
